$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.862.43"
$ws.Range("D3").Value = "1.628.96"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "214.30"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "0.502"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "19.60"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").Value = "1.857.79"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "1.597.34"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "0.543"
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").Value = "62.75"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "25.858.63"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").Value = "192.91"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("D22").Value = "9.93"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "6.25"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "1.79"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").Value = "142.07"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("D28").Value = "6.86"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Value = "0.0500"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("D32").Value = "3.30"
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("D36").Value = "0.900"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("D37").Value = "1.132.47"
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("D38").Value = "0.548"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").Value = "0.997"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.46"
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.802"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "98.94"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("D45").Value = "1.767.07"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "0.0₆0111"
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("D47").Value = "56.00"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("E48").Value = "  +3.60%  "
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("E51").Value = "  +2.30%  "
